# Insert two new rows at 850-851 (pushing existing rows 850:911 down to 852:913),
# then populate the two new rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 850.
$ws.Range("A850:R851").EntireRow.Insert()

# New row 850: Tomate / Larga vida / Primera, Región de Arica y Parinacota
$ws.Range("A850").Value = 11
$ws.Range("B850").Value = "Vega Monumental Concepción"
$ws.Range("C850").Value = "Bíobío"
$ws.Range("D850").Value = 45223
$ws.Range("E850").Value = 8
$ws.Range("F850").Value = 100112020
$ws.Range("G850").Value = "Tomate"
$ws.Range("H850").Value = "Larga vida"
$ws.Range("I850").Value = "Primera"
$ws.Range("J850").Value = 350
$ws.Range("K850").Value = 20000
$ws.Range("L850").Value = 20000
$ws.Range("M850").Value = 20000
$ws.Range("N850").Value = "`$/bandeja 18 kilos"
$ws.Range("O850").Value = "Región de Arica y Parinacota"
$ws.Range("P850").Value = 1111
$ws.Range("Q850").Value = 18
$ws.Range("R850").Value = "Hortaliza"

# New row 851: Tomate / Larga vida / Segunda, Región de Arica y Parinacota
$ws.Range("A851").Value = 11
$ws.Range("B851").Value = "Vega Monumental Concepción"
$ws.Range("C851").Value = "Bíobío"
$ws.Range("D851").Value = 45223
$ws.Range("E851").Value = 8
$ws.Range("F851").Value = 100112020
$ws.Range("G851").Value = "Tomate"
$ws.Range("H851").Value = "Larga vida"
$ws.Range("I851").Value = "Segunda"
$ws.Range("J851").Value = 250
$ws.Range("K851").Value = 18000
$ws.Range("L851").Value = 18000
$ws.Range("M851").Value = 18000
$ws.Range("N851").Value = "`$/bandeja 18 kilos"
$ws.Range("O851").Value = "Región de Arica y Parinacota"
$ws.Range("P851").Value = 1000
$ws.Range("Q851").Value = 18
$ws.Range("R851").Value = "Hortaliza"
